# Daily attendance processing - 2025-11-27 06:37:09
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Recorded By" email lists (values only, same sets, new order) ---
$ws.Range("G2").Value  = "servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value  = "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G4").Value  = "servinaz@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G5").Value  = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G7").Value  = "AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"
$ws.Range("G15").Value = "Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# --- Class statistics shift: one BIOCHEMISTRY LAB/CBL session moved from Pending to Not Recorded ---
$ws.Range("L7").Value = 2    # Missing Sessions 1 -> 2
$ws.Range("L8").Value = 13   # Pending Sessions 14 -> 13

# --- Group statistics mirrors the same shift ---
$ws.Range("P15").Value = 2   # Missing 1 -> 2
$ws.Range("Q15").Value = 13  # Pending 14 -> 13

# --- Row 8 (BIOCHEMISTRY LAB/CBL, session 2) status flips from Pending to Not Recorded ---
$ws.Range("I8").Value = "Not Recorded"

# Re-colour row 8 to match the "Not Recorded" (red) status styling already used
# elsewhere in the sheet (e.g. row 29), reusing that cell's existing format instead
# of fabricating a brand-new style entry.
$ws.Range("A29:I29").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
